# Update average_county_temperature (column AA) values with NOAA data
# Facility 1000032 (KRATON Corporation) - rows 2-13
# Facility 1006856 (LION ELASTOMERS ORANGE PLANT) - rows 38-41
# Facility 1006919 (AMERICAN SYNTHETIC RUBBER) - rows 42-53
# Facility 1013817 (ARLANXEO Orange Site) - rows 70-73

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($r in 2..13) {
    $ws.Range("AA$r").Value = 21.28240740740739
}

foreach ($r in 38..41) {
    $ws.Range("AA$r").Value = 19.65277777777778
}

foreach ($r in 42..53) {
    $ws.Range("AA$r").Value = 13.75752314814816
}

foreach ($r in 70..73) {
    $ws.Range("AA$r").Value = 19.65277777777778
}
